$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 18:18:32"
$ws.Range("E3").Value = "2026-02-12 18:18:35"
$ws.Range("H3").Value = "72%"
$ws.Range("I3").Value = "0.1 mm"
$ws.Range("E4").Value = "2026-02-12 18:18:37"
$ws.Range("J4").Value = "998.5 hPa"
$ws.Range("E5").Value = "2026-02-12 18:18:40"
$ws.Range("H5").Value = "75%"
$ws.Range("E6").Value = "2026-02-12 18:18:42"
$ws.Range("H6").Value = "39%"
$ws.Range("J6").Value = "998.3 hPa"
$ws.Range("E7").Value = "2026-02-12 18:18:45"
$ws.Range("J7").Value = "1001.2 hPa"
$ws.Range("E8").Value = "2026-02-12 18:18:47"
$ws.Range("J8").Value = "1000.5 hPa"
$ws.Range("E9").Value = "2026-02-12 18:18:50"
$ws.Range("E10").Value = "2026-02-12 18:18:53"
$ws.Range("O10").Value = "15.1 °C"
$ws.Range("E11").Value = "2026-02-12 18:18:55"
$ws.Range("H11").Value = "48%"
$ws.Range("E12").Value = "2026-02-12 18:18:58"
$ws.Range("H12").Value = "75%"
$ws.Range("E13").Value = "2026-02-12 18:19:00"
$ws.Range("J13").Value = "1001.1 hPa"
$ws.Range("E14").Value = "2026-02-12 18:19:03"
$ws.Range("E15").Value = "2026-02-12 18:19:05"
$ws.Range("H15").Value = "56%"
$ws.Range("E16").Value = "2026-02-12 18:19:08"
$ws.Range("E17").Value = "2026-02-12 18:19:10"
$ws.Range("O17").Value = "2.2 °C"
$ws.Range("E18").Value = "2026-02-12 18:19:13"
$ws.Range("J18").Value = "998.7 hPa"
$ws.Range("K18").Value = "13.8 MJ/m2"
$ws.Range("E19").Value = "2026-02-12 18:19:15"
$ws.Range("O19").Value = "8.3 °C"
$ws.Range("E20").Value = "2026-02-12 18:19:18"
$ws.Range("E21").Value = "2026-02-12 18:19:20"
$ws.Range("H21").Value = "48%"
$ws.Range("J21").Value = "1001.5 hPa"
$ws.Range("E22").Value = "2026-02-12 18:19:23"
$ws.Range("E23").Value = "2026-02-12 18:19:25"
$ws.Range("E24").Value = "2026-02-12 18:19:28"
$ws.Range("H24").Value = "60%"
$ws.Range("J24").Value = "1006.2 hPa"
$ws.Range("E25").Value = "2026-02-12 18:19:30"
$ws.Range("E26").Value = "2026-02-12 18:19:32"
$ws.Range("J26").Value = "997.7 hPa"
$ws.Range("E27").Value = "2026-02-12 18:19:35"
$ws.Range("E28").Value = "2026-02-12 18:19:38"
$ws.Range("J28").Value = "998.0 hPa"
$ws.Range("O28").Value = "14.4 °C"
$ws.Range("E29").Value = "2026-02-12 18:19:40"
$ws.Range("H29").Value = "54%"
$ws.Range("O29").Value = "15.3 °C"
$ws.Range("E30").Value = "2026-02-12 18:19:42"
$ws.Range("H30").Value = "67%"
$ws.Range("J30").Value = "998.5 hPa"
$ws.Range("O30").Value = "12.6 °C"
$ws.Range("E31").Value = "2026-02-12 18:19:45"
$ws.Range("J31").Value = "998.0 hPa"
$ws.Range("O31").Value = "14.5 °C"
$ws.Range("E32").Value = "2026-02-12 18:19:47"
$ws.Range("K32").Value = "13.9 MJ/m2"
$ws.Range("E33").Value = "2026-02-12 18:19:50"
$ws.Range("J33").Value = "1000.8 hPa"
$ws.Range("O33").Value = "6.7 °C"
$ws.Range("E34").Value = "2026-02-12 18:19:53"
$ws.Range("H34").Value = "58%"
$ws.Range("E35").Value = "2026-02-12 18:19:55"
$ws.Range("J35").Value = "1007.6 hPa"
$ws.Range("E36").Value = "2026-02-12 18:19:58"
$ws.Range("J36").Value = "998.9 hPa"
$ws.Range("E37").Value = "2026-02-12 18:20:00"
$ws.Range("H37").Value = "46%"
$ws.Range("J37").Value = "999.2 hPa"
$ws.Range("E38").Value = "2026-02-12 18:20:02"
$ws.Range("O38").Value = "16.0 °C"
$ws.Range("E39").Value = "2026-02-12 18:20:05"
$ws.Range("H39").Value = "64%"
$ws.Range("E40").Value = "2026-02-12 18:20:07"
$ws.Range("H40").Value = "53%"
$ws.Range("J40").Value = "1002.3 hPa"
$ws.Range("E41").Value = "2026-02-12 18:20:10"
$ws.Range("J41").Value = "1005.2 hPa"
$ws.Range("O41").Value = "17.3 °C"
$ws.Range("E42").Value = "2026-02-12 18:20:12"
$ws.Range("E43").Value = "2026-02-12 18:20:15"
$ws.Range("E44").Value = "2026-02-12 18:20:17"
$ws.Range("I44").Value = "0.5 mm"
$ws.Range("E45").Value = "2026-02-12 18:20:20"
$ws.Range("J45").Value = "1004.3 hPa"
$ws.Range("E46").Value = "2026-02-12 18:20:23"
$ws.Range("J46").Value = "1007.0 hPa"
